# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Estado de Cuenta" detail table (rows 16-25, columns C:F) is
# re-sorted: instead of being grouped by worker (all periods for worker 1,
# then all periods for worker 2), it is now grouped by period (ascending),
# interleaving both workers for each period. The underlying data values
# themselves are unchanged - only the row order/presentation changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New row order: for each period (ascending), worker 1 row then worker 2 row.
# Columns: C = N Doc Trabajador, D = Nombre Trabajador, E = Periodo Mora, F = Valor Mora
$rows = @(
    @{ Row = 16; Doc = "9072343"; Nombre = "ROBERTO PACHECO DE AVILA";     Periodo = "1811"; Mora = 8333  },
    @{ Row = 17; Doc = "9295145"; Nombre = "ANGEL ENRIQUE BARRIOS CUADRO"; Periodo = "1811"; Mora = 8333  },
    @{ Row = 18; Doc = "9072343"; Nombre = "ROBERTO PACHECO DE AVILA";     Periodo = "1812"; Mora = 31249 },
    @{ Row = 19; Doc = "9295145"; Nombre = "ANGEL ENRIQUE BARRIOS CUADRO"; Periodo = "1812"; Mora = 31249 },
    @{ Row = 20; Doc = "9072343"; Nombre = "ROBERTO PACHECO DE AVILA";     Periodo = "1901"; Mora = 31249 },
    @{ Row = 21; Doc = "9295145"; Nombre = "ANGEL ENRIQUE BARRIOS CUADRO"; Periodo = "1901"; Mora = 31249 },
    @{ Row = 22; Doc = "9072343"; Nombre = "ROBERTO PACHECO DE AVILA";     Periodo = "1902"; Mora = 31249 },
    @{ Row = 23; Doc = "9295145"; Nombre = "ANGEL ENRIQUE BARRIOS CUADRO"; Periodo = "1902"; Mora = 31249 },
    @{ Row = 24; Doc = "9072343"; Nombre = "ROBERTO PACHECO DE AVILA";     Periodo = "1903"; Mora = 26041 },
    @{ Row = 25; Doc = "9295145"; Nombre = "ANGEL ENRIQUE BARRIOS CUADRO"; Periodo = "1903"; Mora = 26041 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("C$n").Value = $r.Doc
    $ws.Range("D$n").Value = $r.Nombre
    $ws.Range("E$n").Value = $r.Periodo
    $ws.Range("F$n").Value = $r.Mora
}
